# Add six new "location" rows (251-256) to the LIVE CAM list sheet.
# These are new live-cam entries for border crossings in Chihuahua (Puente
# Internacional Zaragoza, Paso del Norte, Guadalupe, Lerdo) plus a Cozumel
# harbor cam. Columns are: A=Category, B=lat,long, C=Location, D=City,
# E=Country, F=YouTube Link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the six new rows with the same cell formatting (borders / vertical
# alignment) already used by the last existing data row (250), so the
# appended rows visually match the rest of the table.
$ws.Range("A250:F250").Copy($ws.Range("A251:F251"))
$ws.Range("A250:F250").Copy($ws.Range("A252:F252"))
$ws.Range("A250:F250").Copy($ws.Range("A253:F253"))
$ws.Range("A250:F250").Copy($ws.Range("A254:F254"))
$ws.Range("A250:F250").Copy($ws.Range("A255:F255"))
$ws.Range("A250:F250").Copy($ws.Range("A256:F256"))

# Row 251 - Puente Internacional Zaragoza (Sur) BORDER
$ws.Range("E251").Value = "Maxico"
$ws.Range("F251").Value = "GC5RY3zipa4"
$ws.Range("B251").Value = "31.670954411157933, -106.34064206750381"
$ws.Range("A251").Value = "LIVE, TRAFFIC, BORDER"
$ws.Range("C251").Value = "Transmisión en vivo Puente Internacional Zaragoza (Sur) BORDER"
$ws.Range("D251").Value = "Chihuahua"

# Row 252 - Puente Internacional Paso del Norte (Norte)
$ws.Range("E252").Value = "Maxico"
$ws.Range("A252").Value = "LIVE, TRAFFIC, BORDER"
$ws.Range("D252").Value = "Chihuahua"
$ws.Range("F252").Value = "0Pg3S6s76IE"
$ws.Range("C252").Value = "Transmisión en vivo Puente Internacional Paso del Norte (Norte)"
$ws.Range("B252").Value = "31.746326893361665, -106.48663096506293"

# Row 253 - Puente Internacional Paso del Norte (Sur)
$ws.Range("E253").Value = "Maxico"
$ws.Range("A253").Value = "LIVE, TRAFFIC, BORDER"
$ws.Range("D253").Value = "Chihuahua"
$ws.Range("B253").Value = "31.745606549688638, -106.48635172174104"
$ws.Range("C253").Value = "Transmisión en vivo Puente Internacional Paso del Norte (Sur)"
$ws.Range("F253").Value = "IcvugJWPXz8"

# Row 254 - Puente Internacional Guadalupe (Sur)
$ws.Range("E254").Value = "Maxico"
$ws.Range("A254").Value = "LIVE, TRAFFIC, BORDER"
$ws.Range("D254").Value = "Chihuahua"
$ws.Range("B254").Value = "31.431021053324287, -106.15179344486819"
$ws.Range("C254").Value = "Transmisión en vivo Puente Internacional Guadalupe (Sur)"
$ws.Range("F254").Value = "0gT7jvaLCkg"

# Row 255 - Puente Internacional Lerdo (Fila) / Chamizal
$ws.Range("E255").Value = "Maxico"
$ws.Range("A255").Value = "LIVE, TRAFFIC, BORDER"
$ws.Range("F255").Value = "RVXhhbkBGbI"
$ws.Range("D255").Value = "Chamizal"
$ws.Range("C255").Value = "Transmisión en vivo Puente Internacional Lerdo (Fila)"
$ws.Range("B255").Value = "31.74581414672601, -106.48192250193762"

# Row 256 - Cozumel, Quintana Roo En Vivo
$ws.Range("E256").Value = "Maxico"
$ws.Range("F256").Value = "n2masVeEo4A"
$ws.Range("B256").Value = "20.47757950476707, -86.97687618789205"
$ws.Range("D256").Value = "Cozumel"
$ws.Range("C256").Value = "Cozumel, Quintana Roo En Vivo"
$ws.Range("A256").Value = "LIVE, SEA, HARBOR, CRUISE"

# Match the author's final selection state (one row below the new data).
$ws.Range("A257").Select()
